$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.5
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 2.05
$ws.Range("K2").Value = 2.38
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.91
$ws.Range("Z2").Value = 11
$ws.Range("AD2").Value = 8
$ws.Range("AX2").Value = 29

# Row 5 updates
$ws.Range("G5").Value = 3.6
$ws.Range("I5").Value = 2.15
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 2.25
$ws.Range("W5").Value = 8.5
$ws.Range("X5").Value = 17
$ws.Range("AI5").Value = 9
$ws.Range("AJ5").Value = 9.5
$ws.Range("AK5").Value = 19
$ws.Range("AN5").Value = 5.5
$ws.Range("AQ5").Value = 81
$ws.Range("AR5").Value = 126
$ws.Range("AS5").Value = 351
$ws.Range("AT5").Value = 2.25
